$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = 17.59282018238525
$ws.Cells.Item(2,5).Value = 18.3511962890625
$ws.Cells.Item(2,6).Value = 19.34984969049578
$ws.Cells.Item(2,7).Value = 17.00714320229318
$ws.Cells.Item(2,8).Value = 5186523104
$ws.Cells.Item(2,9).Value = "TSM"

$ws.Cells.Item(3,4).Value = 17.91963123757877
$ws.Cells.Item(3,5).Value = 17.12940216064453
$ws.Cells.Item(3,6).Value = 18.07457892669017
$ws.Cells.Item(3,7).Value = 16.52510794633888
$ws.Cells.Item(3,8).Value = 5186523104
$ws.Cells.Item(3,9).Value = "TSM"

$ws.Cells.Item(4,4).Value = 16.23070965230288
$ws.Cells.Item(4,5).Value = 17.01319122314453
$ws.Cells.Item(4,6).Value = 17.6562231200483
$ws.Cells.Item(4,7).Value = 15.65740407005878
$ws.Cells.Item(4,8).Value = 5186523104
$ws.Cells.Item(4,9).Value = "TSM"

$ws.Cells.Item(5,4).Value = 17.23011486182643
$ws.Cells.Item(5,5).Value = 17.31533622741699
$ws.Cells.Item(5,6).Value = 17.36956746817808
$ws.Cells.Item(5,7).Value = 15.84333928489827
$ws.Cells.Item(5,8).Value = 5186523104
$ws.Cells.Item(5,9).Value = "TSM"

$ws.Cells.Item(6,4).Value = 20.05015403506153
$ws.Cells.Item(6,5).Value = 18.27601051330566
$ws.Cells.Item(6,6).Value = 20.30581669025572
$ws.Cells.Item(6,7).Value = 18.14430541955771
$ws.Cells.Item(6,8).Value = 5186523104
$ws.Cells.Item(6,9).Value = "TSM"

$ws.Cells.Item(7,4).Value = 21.40557044963776
$ws.Cells.Item(7,5).Value = 22.32995796203613
$ws.Cells.Item(7,6).Value = 22.90066528088008
$ws.Cells.Item(7,7).Value = 20.69821466199044
$ws.Cells.Item(7,8).Value = 5186523104
$ws.Cells.Item(7,9).Value = "TSM"

$ws.Cells.Item(8,4).Value = 24.62887130900841
$ws.Cells.Item(8,5).Value = 24.99862670898437
$ws.Cells.Item(8,6).Value = 25.41661054090357
$ws.Cells.Item(8,7).Value = 23.72860009134853
$ws.Cells.Item(8,8).Value = 5186523104
$ws.Cells.Item(8,9).Value = "TSM"

$ws.Cells.Item(9,4).Value = 23.51156437329436
$ws.Cells.Item(9,5).Value = 24.84589576721192
$ws.Cells.Item(9,6).Value = 25.27995615309428
$ws.Cells.Item(9,7).Value = 23.29453418035317
$ws.Cells.Item(9,8).Value = 5186523104
$ws.Cells.Item(9,9).Value = "TSM"

$ws.Cells.Item(10,4).Value = 26.47764134629288
$ws.Cells.Item(10,5).Value = 26.58213806152344
$ws.Cells.Item(10,6).Value = 26.79113149198455
$ws.Cells.Item(10,7).Value = 25.31211155723761
$ws.Cells.Item(10,8).Value = 5186523104
$ws.Cells.Item(10,9).Value = "TSM"

$ws.Cells.Item(11,4).Value = 29.08337048277919
$ws.Cells.Item(11,5).Value = 29.84697341918945
$ws.Cells.Item(11,6).Value = 30.38647783407464
$ws.Cells.Item(11,7).Value = 28.47746732683177
$ws.Cells.Item(11,8).Value = 5186523104
$ws.Cells.Item(11,9).Value = "TSM"

$ws.Cells.Item(12,4).Value = 31.3990805691756
$ws.Cells.Item(12,5).Value = 35.13410186767578
$ws.Cells.Item(12,6).Value = 35.40800178313307
$ws.Cells.Item(12,7).Value = 31.34927914535756
$ws.Cells.Item(12,8).Value = 5186523104
$ws.Cells.Item(12,9).Value = "TSM"

$ws.Cells.Item(13,4).Value = 33.62349522680363
$ws.Cells.Item(13,5).Value = 37.60752105712891
$ws.Cells.Item(13,6).Value = 38.65332577954808
$ws.Cells.Item(13,7).Value = 33.54879466332961
$ws.Cells.Item(13,8).Value = 5186523104
$ws.Cells.Item(13,9).Value = "TSM"

$ws.Cells.Item(14,4).Value = 35.97239916639661
$ws.Cells.Item(14,5).Value = 31.91367721557617
$ws.Cells.Item(14,6).Value = 36.75260234352606
$ws.Cells.Item(14,7).Value = 31.5318758629096
$ws.Cells.Item(14,8).Value = 5186523104
$ws.Cells.Item(14,9).Value = "TSM"

$ws.Cells.Item(15,4).Value = 31.15366999532592
$ws.Cells.Item(15,5).Value = 35.45547485351562
$ws.Cells.Item(15,6).Value = 36.06633035569457
$ws.Cells.Item(15,7).Value = 31.00741020551545
$ws.Cells.Item(15,8).Value = 5186523104
$ws.Cells.Item(15,9).Value = "TSM"

$ws.Cells.Item(16,4).Value = 38.32909297949971
$ws.Cells.Item(16,5).Value = 32.77976226806641
$ws.Cells.Item(16,6).Value = 38.68184096024415
$ws.Cells.Item(16,7).Value = 30.94719517108607
$ws.Cells.Item(16,8).Value = 5186523104
$ws.Cells.Item(16,9).Value = "TSM"

$ws.Cells.Item(17,4).Value = 31.14507564582993
$ws.Cells.Item(17,5).Value = 32.36678695678711
$ws.Cells.Item(17,6).Value = 32.91741759169837
$ws.Cells.Item(17,7).Value = 29.43295543851134
$ws.Cells.Item(17,8).Value = 5186523104
$ws.Cells.Item(17,9).Value = "TSM"

$ws.Cells.Item(18,4).Value = 35.3436283757492
$ws.Cells.Item(18,5).Value = 37.70101547241211
$ws.Cells.Item(18,6).Value = 39.26687209341199
$ws.Cells.Item(18,7).Value = 35.05110552995711
$ws.Cells.Item(18,8).Value = 5186523104
$ws.Cells.Item(18,9).Value = "TSM"

$ws.Cells.Item(19,4).Value = 36.90145008144701
$ws.Cells.Item(19,5).Value = 37.86973571777344
$ws.Cells.Item(19,6).Value = 39.86849019185484
$ws.Cells.Item(19,7).Value = 34.85827951902098
$ws.Cells.Item(19,8).Value = 5186523104
$ws.Cells.Item(19,9).Value = "TSM"

$ws.Cells.Item(20,4).Value = 42.01661466299233
$ws.Cells.Item(20,5).Value = 46.19501495361328
$ws.Cells.Item(20,6).Value = 46.57974999364257
$ws.Cells.Item(20,7).Value = 41.3724099322366
$ws.Cells.Item(20,8).Value = 5186523104
$ws.Cells.Item(20,9).Value = "TSM"

$ws.Cells.Item(21,4).Value = 53.69856027779825
$ws.Cells.Item(21,5).Value = 48.5989990234375
$ws.Cells.Item(21,6).Value = 54.63558298385855
$ws.Cells.Item(21,7).Value = 48.29266480941686
$ws.Cells.Item(21,8).Value = 5186523104
$ws.Cells.Item(21,9).Value = "TSM"

$ws.Cells.Item(22,4).Value = 42.59331514957564
$ws.Cells.Item(22,5).Value = 48.32335662841797
$ws.Cells.Item(22,6).Value = 49.74222199352009
$ws.Cells.Item(22,7).Value = 42.11126162334419
$ws.Cells.Item(22,8).Value = 5186523104
$ws.Cells.Item(22,9).Value = "TSM"

$ws.Cells.Item(23,4).Value = 52.2033288817686
$ws.Cells.Item(23,5).Value = 72.28928375244141
$ws.Cells.Item(23,6).Value = 76.42193461658707
$ws.Cells.Item(23,7).Value = 52.01089689692412
$ws.Cells.Item(23,8).Value = 5186523104
$ws.Cells.Item(23,9).Value = "TSM"

$ws.Cells.Item(24,4).Value = 75.58227936684037
$ws.Cells.Item(24,5).Value = 77.24940490722656
$ws.Cells.Item(24,6).Value = 84.0652521243512
$ws.Cells.Item(24,7).Value = 74.24674344576198
$ws.Cells.Item(24,8).Value = 5186523104
$ws.Cells.Item(24,9).Value = "TSM"

$ws.Cells.Item(25,4).Value = 103.103887673482
$ws.Cells.Item(25,5).Value = 112.3996047973633
$ws.Cells.Item(25,6).Value = 125.913093915184
$ws.Cells.Item(25,7).Value = 102.160437374903
$ws.Cells.Item(25,8).Value = 5186523104
$ws.Cells.Item(25,9).Value = "TSM"

$ws.Cells.Item(26,4).Value = 114.020777072652
$ws.Cells.Item(26,5).Value = 108.3851928710938
$ws.Cells.Item(26,6).Value = 118.2822853335467
$ws.Cells.Item(26,7).Value = 105.627747599346
$ws.Cells.Item(26,8).Value = 5186523104
$ws.Cells.Item(26,9).Value = "TSM"

$ws.Cells.Item(27,4).Value = 111.9324772978053
$ws.Cells.Item(27,5).Value = 108.7077789306641
$ws.Cells.Item(27,6).Value = 117.2262034700634
$ws.Cells.Item(27,7).Value = 104.4299299097196
$ws.Cells.Item(27,8).Value = 5186523104
$ws.Cells.Item(27,9).Value = "TSM"

$ws.Cells.Item(28,4).Value = 104.8120521085647
$ws.Cells.Item(28,5).Value = 106.402946472168
$ws.Cells.Item(28,6).Value = 109.631536289722
$ws.Cells.Item(28,7).Value = 100.6757210513928
$ws.Cells.Item(28,8).Value = 5186523104
$ws.Cells.Item(28,9).Value = "TSM"

$ws.Cells.Item(29,4).Value = 116.6440899817867
$ws.Cells.Item(29,5).Value = 115.2345504760742
$ws.Cells.Item(29,6).Value = 136.255485552211
$ws.Cells.Item(29,7).Value = 107.303546399048
$ws.Cells.Item(29,8).Value = 5186523104
$ws.Cells.Item(29,9).Value = "TSM"

$ws.Cells.Item(30,4).Value = 99.87067124890606
$ws.Cells.Item(30,5).Value = 87.746826171875
$ws.Cells.Item(30,6).Value = 100.0406320305404
$ws.Cells.Item(30,7).Value = 85.18797738450533
$ws.Cells.Item(30,8).Value = 5186523104
$ws.Cells.Item(30,9).Value = "TSM"

$ws.Cells.Item(31,4).Value = 75.02463495169746
$ws.Cells.Item(31,5).Value = 83.98506927490234
$ws.Cells.Item(31,6).Value = 84.67797998265208
$ws.Cells.Item(31,7).Value = 69.9938810347539
$ws.Cells.Item(31,8).Value = 5186523104
$ws.Cells.Item(31,9).Value = "TSM"

$ws.Cells.Item(32,4).Value = 64.98689683773776
$ws.Cells.Item(32,5).Value = 58.75357437133789
$ws.Cells.Item(32,6).Value = 71.51612898264993
$ws.Cells.Item(32,7).Value = 56.80625760458374
$ws.Cells.Item(32,8).Value = 5186523104
$ws.Cells.Item(32,9).Value = "TSM"

$ws.Cells.Item(33,4).Value = 72.81426539960387
$ws.Cells.Item(33,5).Value = 89.01868438720703
$ws.Cells.Item(33,6).Value = 94.07776074769215
$ws.Cells.Item(33,7).Value = 70.81751116882097
$ws.Cells.Item(33,8).Value = 5186523104
$ws.Cells.Item(33,9).Value = "TSM"

$ws.Cells.Item(34,4).Value = 89.75228627399213
$ws.Cells.Item(34,5).Value = 81.34735870361328
$ws.Cells.Item(34,6).Value = 90.32161773300452
$ws.Cells.Item(34,7).Value = 78.36558347202754
$ws.Cells.Item(34,8).Value = 5186523104
$ws.Cells.Item(34,9).Value = "TSM"

$ws.Cells.Item(35,4).Value = 99.08191374529088
$ws.Cells.Item(35,5).Value = 96.07796478271484
$ws.Cells.Item(35,6).Value = 103.9754488727992
$ws.Cells.Item(35,7).Value = 93.34533919323496
$ws.Cells.Item(35,8).Value = 5186523104
$ws.Cells.Item(35,9).Value = "TSM"

$ws.Cells.Item(36,4).Value = 84.91729230851597
$ws.Cells.Item(36,5).Value = 84.06986999511719
$ws.Cells.Item(36,6).Value = 92.30055944471192
$ws.Cells.Item(36,7).Value = 82.74516743759578
$ws.Cells.Item(36,8).Value = 5186523104
$ws.Cells.Item(36,9).Value = "TSM"

$ws.Cells.Item(37,4).Value = 100.0655316771504
$ws.Cells.Item(37,5).Value = 110.5467224121094
$ws.Cells.Item(37,6).Value = 116.3108878938039
$ws.Cells.Item(37,7).Value = 96.6892404408629
$ws.Cells.Item(37,8).Value = 5186523104
$ws.Cells.Item(37,9).Value = "TSM"

$ws.Cells.Item(38,4).Value = 134.902511402015
$ws.Cells.Item(38,5).Value = 134.9516448974609
$ws.Cells.Item(38,6).Value = 145.848785480042
$ws.Cells.Item(38,7).Value = 123.5926764458762
$ws.Cells.Item(38,8).Value = 5186523104
$ws.Cells.Item(38,9).Value = "TSM"

$ws.Cells.Item(39,4).Value = 171.7501201123501
$ws.Cells.Item(39,5).Value = 163.4307250976562
$ws.Cells.Item(39,6).Value = 190.7053196752564
$ws.Cells.Item(39,7).Value = 149.9363524534696
$ws.Cells.Item(39,8).Value = 5186523104
$ws.Cells.Item(39,9).Value = "TSM"

$ws.Cells.Item(40,4).Value = 173.442623027762
$ws.Cells.Item(40,5).Value = 188.5103912353516
$ws.Cells.Item(40,6).Value = 210.3354242645952
$ws.Cells.Item(40,7).Value = 168.9608770832319
$ws.Cells.Item(40,8).Value = 5186523104
$ws.Cells.Item(40,9).Value = "TSM"

$ws.Cells.Item(41,4).Value = 195.5205034321339
$ws.Cells.Item(41,5).Value = 207.7479858398437
$ws.Cells.Item(41,6).Value = 224.699699409506
$ws.Cells.Item(41,7).Value = 186.2605711730957
$ws.Cells.Item(41,8).Value = 5186523104
$ws.Cells.Item(41,9).Value = "TSM"

$ws.Cells.Item(42,4).Value = 165.7838078140205
$ws.Cells.Item(42,5).Value = 166.0827178955078
$ws.Cells.Item(42,6).Value = 171.1242808852935
$ws.Cells.Item(42,7).Value = 133.7609007793342
$ws.Cells.Item(42,8).Value = 5186523104
$ws.Cells.Item(42,9).Value = "TSM"

$ws.Cells.Item(43,4).Value = 227.4299926757812
$ws.Cells.Item(43,5).Value = 241.6199951171875
$ws.Cells.Item(43,6).Value = 248.2799987792969
$ws.Cells.Item(43,7).Value = 221.1799926757812
$ws.Cells.Item(43,8).Value = 5186523104
$ws.Cells.Item(43,9).Value = "TSM"
